$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 17:08:56"
$wsZhCn.Range("H2").Value = "2016-03-19 17:09:39"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 17:09:04"
$wsDeDe.Range("H2").Value = "2016-03-19 17:09:53"
